$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Fitness" column (C) values for rows 2-12 per the commit's
# "temp solve of RWheel" fix.
$ws.Range("C2").Value = 4329
$ws.Range("C3").Value = 4471
$ws.Range("C4").Value = 4471
$ws.Range("C5").Value = 4471
$ws.Range("C6").Value = 4471
$ws.Range("C7").Value = 4471
$ws.Range("C8").Value = 4471
$ws.Range("C9").Value = 4471
$ws.Range("C10").Value = 4471
$ws.Range("C11").Value = 4471
$ws.Range("C12").Value = 4471
